$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Remove the NO_LABEL entries from column C on the begin_group rows.
$ws.Range("C3").Clear()
$ws.Range("C8").Clear()
$ws.Range("C12").Clear()
$ws.Range("C16").Clear()
$ws.Range("C21").Clear()

# Move the active selection (bottom-right frozen pane) back to A2.
$ws.Range("A2").Select()
